# Applies the "novedad parte superior" fix: updates the order number, the
# three dates, swaps the product line, and refreshes the payment amounts /
# guarantee note in the lower summary table.

$d = $word.ActiveDocument

# NOTE: Replace must be 1 (wdReplaceOne), not 2 (wdReplaceAll) - with
# wdReplaceAll this runtime's Find.Execute ends up rewriting every matching
# occurrence in the whole document instead of staying inside the supplied
# cell range.
function Replace-InRange($range, [string]$old, [string]$new) {
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1) | Out-Null
}

# --- Table 1: header block (order number + the two dates) ---
$t1 = $d.Tables.Item(1)

Replace-InRange $t1.Cell(1, 1).Range "No.  118" "No.  119"
Replace-InRange $t1.Cell(2, 2).Range "2019-02-19" "2019-02-26"
Replace-InRange $t1.Cell(3, 2).Range "2018-11-27" "2019-02-26"

# --- Table 2: concept / amounts / guarantee block ---
$t2 = $d.Tables.Item(2)

Replace-InRange $t2.Cell(5, 1).Range "Producto 1- (p1)" "Pantalon 3- (p4)"
Replace-InRange $t2.Cell(5, 5).Range "100000" "50000"

Replace-InRange $t2.Cell(6, 3).Range "60000" "39999"

Replace-InRange $t2.Cell(7, 1).Range "2018-11-28" "2019-02-27"
Replace-InRange $t2.Cell(7, 3).Range "40000" "10000"
Replace-InRange $t2.Cell(7, 5).Range "100000" "49999"

Replace-InRange $t2.Cell(8, 2).Range "Ninguna" "10000 - "
